$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 5408
$ws.Range("F5").Value = 7581
$ws.Range("F12").Value = 4376
$ws.Range("F13").Value = 1782
$ws.Range("F14").Value = 113
$ws.Range("F16").Value = 2965
$ws.Range("F19").Value = 217
$ws.Range("F20").Value = 535
$ws.Range("F22").Value = 472
$ws.Range("F23").Value = 329
$ws.Range("F24").Value = 112
$ws.Range("F26").Value = 1224
$ws.Range("F28").Value = 1416
$ws.Range("F31").Value = 33
$ws.Range("F34").Value = 18
$ws.Range("F38").Value = 3027
$ws.Range("F40").Value = 39
$ws.Range("F41").Value = 125
$ws.Range("F42").Value = 46
$ws.Range("F43").Value = 460

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 22

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5408
$ws.Range("F5").Value = 7581
$ws.Range("F12").Value = 4376
$ws.Range("F13").Value = 1782
$ws.Range("F14").Value = 113
$ws.Range("F16").Value = 2965
$ws.Range("F19").Value = 217
$ws.Range("F20").Value = 535
$ws.Range("F22").Value = 472
$ws.Range("F24").Value = 329
$ws.Range("F25").Value = 112
$ws.Range("F27").Value = 1224
$ws.Range("F29").Value = 1416
$ws.Range("F32").Value = 33
$ws.Range("F35").Value = 18
$ws.Range("F39").Value = 3027
$ws.Range("F40").Value = 22
$ws.Range("F42").Value = 39
$ws.Range("F43").Value = 125
$ws.Range("F44").Value = 46
$ws.Range("F45").Value = 460
